$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect to allow cell edits, then
# restore protection afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A80).
$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-03-29 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for each holding row.
$ws.Range("D2").Value = 0.07423426931676005
$ws.Range("E2").Value = 0.001485025987954858
$ws.Range("D3").Value = 0.04486069802481601
$ws.Range("E3").Value = 0.007765323407699087
$ws.Range("D4").Value = 0.03591799490133487
$ws.Range("E4").Value = -0.005243572395128515
$ws.Range("D5").Value = 0.03347841220167635
$ws.Range("E5").Value = 0.004589519037137668
$ws.Range("D6").Value = 0.03191457462520421
$ws.Range("E6").Value = -0.01553936424011859
$ws.Range("D7").Value = 0.02976078253221159
$ws.Range("E7").Value = 0.01040138685158021
$ws.Range("D8").Value = 0.03070713013511776
$ws.Range("E8").Value = 0.006730127933062491
$ws.Range("D9").Value = 0.02893530676654949
$ws.Range("E9").Value = 0.007983798730381464
$ws.Range("D10").Value = 0.02631370587663488
$ws.Range("E10").Value = -0.001139528994682171
$ws.Range("D11").Value = 0.02787303586939139
$ws.Range("E11").Value = 0.01565153193060165
$ws.Range("D12").Value = 0.02368402073331724
$ws.Range("E12").Value = -0.005755163511187766
$ws.Range("D13").Value = 0.02387882674259267
$ws.Range("E13").Value = -0.009565667011375401
$ws.Range("D14").Value = 0.01989843433979185
$ws.Range("E14").Value = -0.008480085096324252
$ws.Range("D15").Value = 0.01894282659207857
$ws.Range("E15").Value = -0.01601556033769236
$ws.Range("D16").Value = 0.02137824467634765
$ws.Range("E16").Value = -0.0009511128019784598
$ws.Range("D17").Value = 0.01853812396717538
$ws.Range("E17").Value = 0.01139643306445648
$ws.Range("D18").Value = 0.0181783011975278
$ws.Range("E18").Value = 0.01478626489138057
$ws.Range("D19").Value = 0.01555140879629485
$ws.Range("E19").Value = -0.005371686016288346
$ws.Range("D20").Value = 0.01401985923795965
$ws.Range("E20").Value = 0.007076806944706471
$ws.Range("D21").Value = 0.01525336932083648
$ws.Range("E21").Value = 0.02755988976044099
$ws.Range("D22").Value = 0.01435136262295904
$ws.Range("E22").Value = -0.008821762167476765
$ws.Range("D23").Value = 0.01315293330030385
$ws.Range("E23").Value = -0.01267265656430183
$ws.Range("D24").Value = 0.01514753909446955
$ws.Range("E24").Value = 0.008248102936324608
$ws.Range("D25").Value = 0.01475665319356427
$ws.Range("E25").Value = 0.006852953676158124
$ws.Range("D26").Value = 0.01207111320855296
$ws.Range("E26").Value = -0.01591089896579145
$ws.Range("D27").Value = 0.01215577738964651
$ws.Range("E27").Value = 0.0001612253123739471
$ws.Range("D28").Value = 0.01239551224965272
$ws.Range("E28").Value = -0.005703715122552966
$ws.Range("D29").Value = 0.0118957094074448
$ws.Range("E29").Value = 0.0001812250815511618
$ws.Range("D30").Value = 0.01144985058339892
$ws.Range("E30").Value = -0.02386046591239765
$ws.Range("D31").Value = 0.01276136146274619
$ws.Range("E31").Value = -0.0002303616678184239
$ws.Range("D32").Value = 0.01339678378022433
$ws.Range("E32").Value = -0.00529936473453807
$ws.Range("D33").Value = 0.01127042915332683
$ws.Range("E33").Value = -0.003238707994609569
$ws.Range("D34").Value = 0.01168728265607215
$ws.Range("E34").Value = 0.01117222413200403
$ws.Range("D35").Value = 0.0090941971327064
$ws.Range("E35").Value = -0.01199269447721885
$ws.Range("D36").Value = 0.01117449601294421
$ws.Range("E36").Value = 0.01527149321266985
$ws.Range("D37").Value = 0.01115186010341573
$ws.Range("E37").Value = -0.003646588462721367
$ws.Range("D38").Value = 0.01019277367696538
$ws.Range("E38").Value = -0.001009445525993136
$ws.Range("D39").Value = 0.009447111540355026
$ws.Range("E39").Value = -0.01954194672641274
$ws.Range("D40").Value = 0.009569404246379037
$ws.Range("E40").Value = 0.001310723355451904
$ws.Range("D41").Value = 0.00912227153997874
$ws.Range("E41").Value = -0.01293330325588382
$ws.Range("D42").Value = 0.00912227153997874
$ws.Range("E42").Value = 0.003759681179035912
$ws.Range("D43").Value = 0.009858575540831653
$ws.Range("E43").Value = 0.01175862514536763
$ws.Range("D44").Value = 0.009637409965914827
$ws.Range("E44").Value = 0.01252669039145893
$ws.Range("D45").Value = 0.009478860608265105
$ws.Range("E45").Value = -0.005210271678451894
$ws.Range("D46").Value = 0.009590864264503444
$ws.Range("E46").Value = 0.01020689655172413
$ws.Range("D47").Value = 0.009040253114544365
$ws.Range("E47").Value = -0.03277835587929245
$ws.Range("D48").Value = 0.007263481202983942
$ws.Range("E48").Value = 0.005558253736981378
$ws.Range("D49").Value = 0.00830522299604587
$ws.Range("E49").Value = -0.004341926729986678
$ws.Range("D50").Value = 0.008076904081754242
$ws.Range("E50").Value = -0.007006369426751702
$ws.Range("D51").Value = 0.007871123086040757
$ws.Range("E51").Value = 0.003137254901960818
$ws.Range("D52").Value = 0.007815562217198117
$ws.Range("E52").Value = 0.004476040021063721
$ws.Range("D53").Value = 0.007005373039789028
$ws.Range("E53").Value = -0.007133864876206486
$ws.Range("D54").Value = 0.007416935031215996
$ws.Range("E54").Value = 0.001268331351565566
$ws.Range("D55").Value = 0.006741532405999241
$ws.Range("E55").Value = -0.004375885751662589
$ws.Range("D56").Value = 0.006584207935228765
$ws.Range("E56").Value = 0.0009524943445649026
$ws.Range("D57").Value = 0.006883668279467054
$ws.Range("E57").Value = -0.007174581482746967
$ws.Range("D58").Value = 0.006439769274427967
$ws.Range("E58").Value = -0.01972062448644196
$ws.Range("D59").Value = 0.005689697545052319
$ws.Range("E59").Value = -0.01269299990527617
$ws.Range("D60").Value = 0.006603365166020187
$ws.Range("E60").Value = -0.02203672787979949
$ws.Range("D61").Value = 0.005520810142141753
$ws.Range("E61").Value = -0.02218672346467876
$ws.Range("D62").Value = 0.005835557074633043
$ws.Range("E62").Value = 0.01356797420741551
$ws.Range("D63").Value = 0.005374999608036198
$ws.Range("E63").Value = 0.01298038357762699
$ws.Range("D64").Value = 0.005012433091779108
$ws.Range("E64").Value = 0.00203315608382848
$ws.Range("D65").Value = 0.004803908349456111
$ws.Range("E65").Value = 0.008240861618798778
$ws.Range("D66").Value = 0.004525418068590529
$ws.Range("E66").Value = 0.004352345069507768
$ws.Range("D67").Value = 0.004467211444088715
$ws.Range("E67").Value = 0.005527770465912241
$ws.Range("D68").Value = 0.003542813814058644
$ws.Range("E68").Value = -0.02223789569763102
$ws.Range("D69").Value = 0.004135316095287968
$ws.Range("E69").Value = -0.00753536646051034
$ws.Range("D70").Value = 0.003705821759948825
$ws.Range("E70").Value = -0.01599767306573596
$ws.Range("D71").Value = 0.003073339170980881
$ws.Range("E71").Value = 0.008034817542685202
$ws.Range("D72").Value = 0.002516260618299379
$ws.Range("E72").Value = 0.008489592460618534
$ws.Range("D73").Value = 0.002489215116005607
$ws.Range("E73").Value = 0.01161303021356175
$ws.Range("D74").Value = 0.002298328744743765
$ws.Range("E74").Value = 0.0004903110277345579
$ws.Range("D75").Value = 0.001922190481864619
$ws.Range("E75").Value = -0.028242251223491
$ws.Range("D76").Value = 0.001787550916097796
$ws.Range("E76").Value = -0.01650038372985418
$ws.Range("D77").Value = 0.9999999999999999
$ws.Range("E77").Value = -0.0002016212798716666

# Restore sheet protection.
$ws.Protect()
